$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("G2").Value = 0.1063306666666667
$ws.Range("H2").Value = 0.318992
$ws.Range("I2").Value = 0.2215797714128531
$ws.Range("J2").Value = 0.2215797714128531
$ws.Range("M2").Value = 7.236245333333333
$ws.Range("N2").Value = 21.708736
$ws.Range("O2").Value = 0.6630432242917509
$ws.Range("P2").Value = 0.6630432242917509
$ws.Range("Q2").Value = 0.7694347904568889
$ws.Range("R2").Value = 6.924913114111999
$ws.Range("S2").Value = 0.1469169660754072
$ws.Range("T2").Value = 0.1469169660754072
$ws.Range("G3").Value = 0.1063306666666667
$ws.Range("H3").Value = 0.318992
$ws.Range("I3").Value = 0.2215797714128531
$ws.Range("J3").Value = 0.2215797714128531
$ws.Range("O3").Value = 0.003787629702975075
$ws.Range("P3").Value = 0.003787629702975075
$ws.Range("Q3").Value = 0.004395390768
$ws.Range("R3").Value = 0.039558516912
$ws.Range("S3").Value = 0.0008392621237817497
$ws.Range("T3").Value = 0.0008392621237817498
$ws.Range("G4").Value = 0.1063306666666667
$ws.Range("H4").Value = 0.318992
$ws.Range("I4").Value = 0.2215797714128531
$ws.Range("J4").Value = 0.2215797714128531
$ws.Range("M4").Value = 3.636103333333333
$ws.Range("N4").Value = 10.90831
$ws.Range("O4").Value = 0.3331691460052741
$ws.Range("P4").Value = 0.3331691460052741
$ws.Range("Q4").Value = 0.3866292915022222
$ws.Range("R4").Value = 3.47966362352
$ws.Range("S4").Value = 0.07382354321366411
$ws.Range("T4").Value = 0.07382354321366411
$ws.Range("I5").Value = 0.6204562851740659
$ws.Range("J5").Value = 0.620456285174066
$ws.Range("M5").Value = 7.236245333333333
$ws.Range("N5").Value = 21.708736
$ws.Range("O5").Value = 0.6630432242917509
$ws.Range("P5").Value = 0.6630432242917509
$ws.Range("Q5").Value = 2.154531745955555
$ws.Range("R5").Value = 19.3907857136
$ws.Range("S5").Value = 0.4113893358538948
$ws.Range("T5").Value = 0.4113893358538948
$ws.Range("I6").Value = 0.6204562851740659
$ws.Range("J6").Value = 0.620456285174066
$ws.Range("O6").Value = 0.003787629702975075
$ws.Range("P6").Value = 0.003787629702975075
$ws.Range("S6").Value = 0.002350058655122866
$ws.Range("T6").Value = 0.002350058655122866
$ws.Range("I7").Value = 0.6204562851740659
$ws.Range("J7").Value = 0.620456285174066
$ws.Range("M7").Value = 3.636103333333333
$ws.Range("N7").Value = 10.90831
$ws.Range("O7").Value = 0.3331691460052741
$ws.Range("P7").Value = 0.3331691460052741
$ws.Range("Q7").Value = 1.082619466638889
$ws.Range("R7").Value = 9.74357519975
$ws.Range("S7").Value = 0.2067168906650483
$ws.Range("T7").Value = 0.2067168906650484
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.075803
$ws.Range("H8").Value = 0.227409
$ws.Range("I8").Value = 0.1579639434130809
$ws.Range("J8").Value = 0.1579639434130809
$ws.Range("M8").Value = 7.236245333333333
$ws.Range("N8").Value = 21.708736
$ws.Range("O8").Value = 0.6630432242917509
$ws.Range("P8").Value = 0.6630432242917509
$ws.Range("Q8").Value = 0.5485291050026666
$ws.Range("R8").Value = 4.936761945023999
$ws.Range("S8").Value = 0.1047369223624488
$ws.Range("T8").Value = 0.1047369223624489
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.075803
$ws.Range("H9").Value = 0.227409
$ws.Range("I9").Value = 0.1579639434130809
$ws.Range("J9").Value = 0.1579639434130809
$ws.Range("O9").Value = 0.003787629702975075
$ws.Range("P9").Value = 0.003787629702975075
$ws.Range("Q9").Value = 0.003133468611
$ws.Range("R9").Value = 0.028201217499
$ws.Range("S9").Value = 0.0005983089240704591
$ws.Range("T9").Value = 0.0005983089240704592
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.075803
$ws.Range("H10").Value = 0.227409
$ws.Range("I10").Value = 0.1579639434130809
$ws.Range("J10").Value = 0.1579639434130809
$ws.Range("M10").Value = 3.636103333333333
$ws.Range("N10").Value = 10.90831
$ws.Range("O10").Value = 0.3331691460052741
$ws.Range("P10").Value = 0.3331691460052741
$ws.Range("Q10").Value = 0.2756275409766666
$ws.Range("R10").Value = 2.48064786879
$ws.Range("S10").Value = 0.0526287121265616
$ws.Range("T10").Value = 0.05262871212656161
